$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.925.91'
$ws.Range("E2").Value = '  +1.12%  '
$ws.Range("D3").Value = '3.517.88'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.95'
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '181.25'
$ws.Range("E6").Value = '  +4.73%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.518.05'
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.594'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.142'
$ws.Range("E10").Value = '  +7.15%  '
$ws.Range("E11").Value = '  -1.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.440'
$ws.Range("E12").Value = '  +0.93%  '
$ws.Range("D13").Value = '4.121.11'
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.47'
$ws.Range("E14").Value = '  +11.49%  '
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("D16").Value = '67.870.85'
$ws.Range("E16").Value = '  +1.12%  '
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").Value = '3.510.84'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.36'
$ws.Range("E19").Value = '  +0.35%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.52'
$ws.Range("E20").Value = '  +2.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '401.01'
$ws.Range("E21").Value = '  +1.17%  '
$ws.Range("E22").Value = '  -0.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.84'
$ws.Range("E23").Value = '  +0.96%  '
$ws.Range("E24").Value = '  +1.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  +0.74%  '
$ws.Range("E27").Value = '  +0.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.47'
$ws.Range("E28").Value = '  +2.34%  '
$ws.Range("E29").Value = '  -2.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.996'
$ws.Range("E30").Value = '  -0.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.28'
$ws.Range("E31").Value = '  -0.85%  '
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.10'
$ws.Range("E33").Value = '  +1.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.95'
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("E35").Value = '  +1.66%  '
$ws.Range("E36").Value = '  +0.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.63'
$ws.Range("E37").Value = '  -2.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '163.15'
$ws.Range("E38").Value = '  -0.29%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.883'
$ws.Range("E39").Value = '  -0.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.93'
$ws.Range("E40").Value = '  +0.84%  '
$ws.Range("E41").Value = '  +8.65%  '
$ws.Range("E42").Value = '  -1.50%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.71'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.889.08'
$ws.Range("E44").Value = '  +2.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.41'
$ws.Range("E45").Value = '  -0.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0737'
$ws.Range("E46").Value = '  -1.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.98'
$ws.Range("E47").Value = '  -0.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '42.39'
$ws.Range("E48").Value = '  -1.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '348.99'
$ws.Range("E49").Value = '  +2.75%  '
$ws.Range("E50").Value = '  -0.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.07'
$ws.Range("E51").Value = '  -1.16%  '
